$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows after row 2 (rows 3 and 4), shifting existing content down.
$ws.Rows("3:4").Insert()

# Row 2 becomes the "Albatros / damas" entry for Ferrante Keller, María Estefanía
$ws.Range("B2").Value = "Albatros"
$ws.Range("C2").Value = "damas"
$ws.Range("D2").Value = 1
$ws.Range("E2").Value = "Ferrante Keller, María Estefanía"
$ws.Range("F2").Value = 115
$ws.Range("G2").Value = ""
$ws.Range("H2").Value = 115

# Row 3: new "Prejuveniles / damas" entry for the same player
$ws.Range("A3").Value = "####1er Torneo Federativo - C.A.E. - Sub 23, Prejuveniles y sub 23 (28 de Febrero y 1 de Marzo) - Juniors (Domingo 1 de Marzo)"
$ws.Range("B3").Value = "Prejuveniles"
$ws.Range("C3").Value = "damas"
$ws.Range("D3").Value = 1
$ws.Range("E3").Value = "Ferrante Keller, María Estefanía"
$ws.Range("F3").Value = ""
$ws.Range("G3").Value = 115
$ws.Range("H3").Value = 115

# Row 4: the original "Prejuveniles / caballeros" entry, moved down
$ws.Range("A4").Value = "####1er Torneo Federativo - C.A.E. - Sub 23, Prejuveniles y sub 23 (28 de Febrero y 1 de Marzo) - Juniors (Domingo 1 de Marzo)"
$ws.Range("B4").Value = "Prejuveniles"
$ws.Range("C4").Value = "caballeros"
$ws.Range("D4").Value = 1
$ws.Range("E4").Value = "Kern Pascuali, Juan Daniel"
$ws.Range("F4").Value = 92
$ws.Range("G4").Value = ""
$ws.Range("H4").Value = 92
